$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")
$hungary  = $wb.Worksheets.Item("Hungary")

# --- Norway: duplicate the "Slovakia" template sheet and place it right
# --- after "Hungary" (i.e. at the end of the tab strip) ---------------
$slovakia.Copy($null, $hungary) | Out-Null
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"

# Norway's repeater list has one extra product ("MZXSDR240") that sits
# between "MZX64DR" (row 14) and "MZXDR240" (row 15) - insert a row and
# clone the formatting of the row above it before writing the new value.
$norway.Rows.Item(15).Insert() | Out-Null
$norway.Range("A14").Copy() | Out-Null
$norway.Range("A15").PasteSpecial(-4122) | Out-Null
$norway.Range("A15").Value = "MZXSDR240"

# Fill in the market-specific cells (B4 before B2 so new shared strings
# land in the same order as the target workbook).
$norway.Range("B4").Value = "NGC-2931/T3061/T3060"
$norway.Range("B2").Value = "Norway Market"

# --- Poland: duplicate the same template sheet and place it after Norway
$slovakia2 = $wb.Worksheets.Item("Slovakia")
$slovakia2.Copy($null, $norway) | Out-Null
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/T3104/T3103"
$poland.Range("B2").Value = "Poland Market"

# Restore the expected selection/active-sheet state: Norway is the
# selected tab, both new sheets keep a B11 selection.
$norway.Activate()
$norway.Range("B11").Select() | Out-Null
$poland.Range("B11").Select() | Out-Null
$norway.Activate()
$norway.Range("B11").Select() | Out-Null
